# New PO forecast model
# Update "Weekly Quantity", "Monthly Trend" and "PO Forecast" sheets with
# refreshed forecast data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append 4 new weekly rows (6-9)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @(45662.99999999999, 15),
    @(45669.99999999999, 15),
    @(45676.99999999999, 6),
    @(45683.99999999999, 11)
)

$r = 6
foreach ($row in $weeklyNewRows) {
    $wsWeekly.Range("A5").Copy($wsWeekly.Range("A$r"))
    $wsWeekly.Range("A$r").Value = $row[0]
    $wsWeekly.Range("B$r").Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append 1 new monthly row (6)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Range("A5").Copy($wsMonthly.Range("A6"))
$wsMonthly.Range("A6").Value = 45688.99999999999
$wsMonthly.Range("B6").Value = 47

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" - refresh the forecast curve
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Updated values for the existing date rows 2-5 (dates unchanged)
$wsForecast.Range("B2").Value = 97
$wsForecast.Range("B3").Value = 66
$wsForecast.Range("B4").Value = 41
$wsForecast.Range("B5").Value = 38

# Rows 6-13 now carry shifted dates/values (one week further out each)
$forecastRows = @(
    @(6,  45662.99999999999, 24),
    @(7,  45669.99999999999, 20),
    @(8,  45676.99999999999, 17),
    @(9,  45683.99999999999, 13),
    @(10, 45690.99999999999, 10),
    @(11, 45697.99999999999, 6),
    @(12, 45704.99999999999, 3),
    @(13, 45711.99999999999, 0),
    @(14, 45718.99999999999, 0),
    @(15, 45725.99999999999, 0),
    @(16, 45732.99999999999, 0),
    @(17, 45739.99999999999, 0)
)

foreach ($row in $forecastRows) {
    $rowNum = $row[0]
    if ($rowNum -gt 13) {
        $wsForecast.Range("A5").Copy($wsForecast.Range("A$rowNum"))
    }
    $wsForecast.Range("A$rowNum").Value = $row[1]
    $wsForecast.Range("B$rowNum").Value = $row[2]
}
